$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2. Muebles -> Sillas x7 becomes Sillas x10, purchase cost 35000 -> 50000
$ws.Range("B12").Value2 = "   Sillas x10"
$ws.Range("C12").Value2 = 50000

# 4. Gastos del profesional -> "Seguro anti incendios" becomes "Seguro (robo e incendios)", monthly cost 900 -> 1500
$ws.Range("B21").Value2 = "Seguro (robo e incendios)"
$ws.Range("D21").Value2 = 1500

# Cursos / Certificaciones x7 monthly cost formula 7*3000 -> 7*10000
$ws.Range("D23").Formula = "=7*10000"
